$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.279.86'
$ws.Range('E2').Value = '  +3.02%  '

$ws.Range('D3').Value = '2.305.89'
$ws.Range('E3').Value = '  +2.76%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.77'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.94%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.57'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +8.59%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.526'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.78%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +6.73%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.85'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +3.50%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.50'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.69%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0809'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.68%  '

$ws.Range('E13').Value = '  -0.59%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.97'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.13%  '

$ws.Range('D15').Value = '2.662.55'
$ws.Range('E15').Value = '  +2.73%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.13'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +5.06%  '

$ws.Range('D17').Value = '2.301.22'
$ws.Range('E17').Value = '  +2.71%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.800'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.19%  '

$ws.Range('D19').Value = '43.232.41'
$ws.Range('E19').Value = '  +3.28%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.00'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.87%  '

$ws.Range('D21').Value = '0.0₃0921'
$ws.Range('E21').Value = '  +2.67%  '

$ws.Range('E22').Value = '  +4.56%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.79'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.08'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.16%  '

$ws.Range('E25').Value = '  +3.47%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.60'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.49%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.11%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.72'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +6.53%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.32'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +9.71%  '

$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.31'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.22%  '

$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.55'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.17%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.86'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.98%  '

$ws.Range('E33').Value = '  +1.99%  '

$ws.Range('E34').Value = '  -0.03%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.23'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +4.57%  '

$ws.Range('E36').Value = '  +6.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0735'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.57%  '

$ws.Range('E38').Value = '  -2.22%  '

$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.49'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +11.34%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.106'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.81%  '

$ws.Range('E41').Value = '  +4.03%  '

$ws.Range('E42').Value = '  +0.73%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.46'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +13.04%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0289'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.51%  '

$ws.Range('D45').Value = '1.962.55'
$ws.Range('E45').Value = '  +1.26%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.72'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.46%  '

$ws.Range('E47').Value = '  +6.54%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.16'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.47%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '57.09'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +7.03%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.94'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.08%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.58'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +8.41%  '
